$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 -> Emails ---
$wsEmails = $wb.Worksheets.Item("Sheet2")
$wsEmails.Name = "Emails"

$ws = $wb.Worksheets.Item("Users")

# --- Append new user rows 53-80 ---
$ws.Cells.Item(53,1).Value = 'SearchOpenWebUser1'
$ws.Cells.Item(53,2).Value = 'Password1'
$ws.Cells.Item(53,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(53,6).Value = 'N'
$ws.Cells.Item(53,7).Value = 'SearchOpenWeb@mailinator.com '

$ws.Cells.Item(54,1).Value = 'FFHUser1'
$ws.Cells.Item(54,2).Value = 'Password1'
$ws.Cells.Item(54,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(54,6).Value = 'N'
$ws.Cells.Item(54,7).Value = 'FFHUser1@mailinator.com '

$ws.Cells.Item(55,1).Value = 'FFHUser2'
$ws.Cells.Item(55,2).Value = 'Password1'
$ws.Cells.Item(55,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(55,6).Value = 'N'
$ws.Cells.Item(55,7).Value = 'FFHUser2@mailinator.com'

$ws.Cells.Item(56,1).Value = 'FFHUser3'
$ws.Cells.Item(56,2).Value = 'Password1'
$ws.Cells.Item(56,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(56,6).Value = 'N'
$ws.Cells.Item(56,7).Value = 'FFHUser3@mailinator.com'

$ws.Cells.Item(57,1).Value = 'FFHUser4'
$ws.Cells.Item(57,2).Value = 'Password1'
$ws.Cells.Item(57,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(57,6).Value = 'N'
$ws.Cells.Item(57,7).Value = 'FFHUser4@mailinator.com'

$ws.Cells.Item(58,1).Value = 'FrontEndUser1'
$ws.Cells.Item(58,2).Value = 'Password1'
$ws.Cells.Item(58,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(58,6).Value = 'N'
$ws.Cells.Item(58,7).Value = 'FrontEndUser1@mailinator.com'

$ws.Cells.Item(59,1).Value = 'FrontEndUser2'
$ws.Cells.Item(59,2).Value = 'Password1'
$ws.Cells.Item(59,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(59,6).Value = 'N'
$ws.Cells.Item(59,7).Value = 'FrontEndUser2@mailinator.com'

$ws.Cells.Item(60,1).Value = 'FrontEndUser3'
$ws.Cells.Item(60,2).Value = 'Password1'
$ws.Cells.Item(60,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(60,6).Value = 'N'
$ws.Cells.Item(60,7).Value = 'FrontEndUser3@mailinator.com'

$ws.Cells.Item(61,1).Value = 'FrontEndUser4'
$ws.Cells.Item(61,2).Value = 'Password1'
$ws.Cells.Item(61,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(61,6).Value = 'N'
$ws.Cells.Item(61,7).Value = 'FrontEndUser4@mailinator.com'

$ws.Cells.Item(62,1).Value = 'FrontEndUser5'
$ws.Cells.Item(62,2).Value = 'Password1'
$ws.Cells.Item(62,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(62,6).Value = 'N'
$ws.Cells.Item(62,7).Value = 'FrontEndUser5@mailinator.com'

$ws.Cells.Item(63,1).Value = 'FrontEndUser6'
$ws.Cells.Item(63,2).Value = 'Password1'
$ws.Cells.Item(63,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(63,6).Value = 'N'
$ws.Cells.Item(63,7).Value = 'FrontEndUser6@mailinator.com'

$ws.Cells.Item(64,1).Value = 'FrontEndUser7'
$ws.Cells.Item(64,2).Value = 'Password1'
$ws.Cells.Item(64,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(64,6).Value = 'N'
$ws.Cells.Item(64,7).Value = 'FrontEndUser7@mailinator.com'

$ws.Cells.Item(65,1).Value = 'FrontEndUser8'
$ws.Cells.Item(65,2).Value = 'Password1'
$ws.Cells.Item(65,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(65,6).Value = 'N'
$ws.Cells.Item(65,7).Value = 'FrontEndUser8@mailinator.com'

$ws.Cells.Item(66,1).Value = 'FrontEndUser9'
$ws.Cells.Item(66,2).Value = 'Password1'
$ws.Cells.Item(66,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(66,6).Value = 'N'
$ws.Cells.Item(66,7).Value = 'FrontEndUser9@mailinator.com'

$ws.Cells.Item(67,1).Value = 'FrontEndUser10'
$ws.Cells.Item(67,2).Value = 'Password1'
$ws.Cells.Item(67,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(67,6).Value = 'N'
$ws.Cells.Item(67,7).Value = 'FrontEndUser10@mailinator.com'

$ws.Cells.Item(68,1).Value = 'UrlUser1'
$ws.Cells.Item(68,2).Value = 'Password1'
$ws.Cells.Item(68,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(68,6).Value = 'N'
$ws.Cells.Item(68,7).Value = 'UrlUser1@mailinator.com'

$ws.Cells.Item(69,1).Value = 'UrlUser2'
$ws.Cells.Item(69,2).Value = 'Password1'
$ws.Cells.Item(69,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(69,6).Value = 'N'
$ws.Cells.Item(69,7).Value = 'UrlUser2@mailinator.com'

$ws.Cells.Item(70,1).Value = 'UrlUser3'
$ws.Cells.Item(70,2).Value = 'Password1'
$ws.Cells.Item(70,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(70,6).Value = 'N'
$ws.Cells.Item(70,7).Value = 'UrlUser3@mailinator.com'

$ws.Cells.Item(71,1).Value = 'LinkingUser1'
$ws.Cells.Item(71,2).Value = 'Password1'
$ws.Cells.Item(71,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(71,6).Value = 'N'
$ws.Cells.Item(71,7).Value = 'LinkingUser1@mailinator.com '

$ws.Cells.Item(72,1).Value = 'LoginUser1'
$ws.Cells.Item(72,2).Value = 'Password1'
$ws.Cells.Item(72,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(72,6).Value = 'N'
$ws.Cells.Item(72,7).Value = 'LoginUser1@mailinator.com '

$ws.Cells.Item(73,1).Value = 'LoginUser2'
$ws.Cells.Item(73,2).Value = 'Password1'
$ws.Cells.Item(73,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(73,6).Value = 'N'
$ws.Cells.Item(73,7).Value = 'LoginUser2@mailinator.com '

$ws.Cells.Item(74,1).Value = 'LoginUser3'
$ws.Cells.Item(74,2).Value = 'Password1'
$ws.Cells.Item(74,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(74,6).Value = 'N'
$ws.Cells.Item(74,7).Value = 'LoginUser3@mailinator.com '

$ws.Cells.Item(75,1).Value = 'LoginUser4'
$ws.Cells.Item(75,2).Value = 'Password1'
$ws.Cells.Item(75,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(75,6).Value = 'N'
$ws.Cells.Item(75,7).Value = 'LoginUser4@mailinator.com '

$ws.Cells.Item(76,1).Value = 'LoginUser5'
$ws.Cells.Item(76,2).Value = 'Password1'
$ws.Cells.Item(76,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(76,6).Value = 'N'
$ws.Cells.Item(76,7).Value = 'LoginUser5@mailinator.com '

$ws.Cells.Item(77,1).Value = 'LoginUser6'
$ws.Cells.Item(77,2).Value = 'Password1'
$ws.Cells.Item(77,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(77,6).Value = 'N'
$ws.Cells.Item(77,7).Value = 'LoginUser6@mailinator.com '

$ws.Cells.Item(78,1).Value = 'LoginUser7'
$ws.Cells.Item(78,2).Value = 'Password1'
$ws.Cells.Item(78,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(78,6).Value = 'N'
$ws.Cells.Item(78,7).Value = 'LoginUser7@mailinator.com '

$ws.Cells.Item(79,1).Value = 'CpetUser1'
$ws.Cells.Item(79,2).Value = 'Password1'
$ws.Cells.Item(79,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(79,6).Value = 'N'
$ws.Cells.Item(79,7).Value = 'CpetUser1@mailinator.com '

$ws.Cells.Item(80,1).Value = 'CpetUser2'
$ws.Cells.Item(80,2).Value = 'Password1'
$ws.Cells.Item(80,5).Value = 'THIS IS IN USE 24/7 - DO NOT USE!'
$ws.Cells.Item(80,6).Value = 'N'
$ws.Cells.Item(80,7).Value = 'CpetUser2@mailinator.com '
